# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-02-13 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-14 Wednesday", 2)

$tbl = $d.Tables.Item(1)

# Row 1 of the table holds the first line of problems.
$tbl.Rows.Item(1).Cells.Item(1).Range.Text = "74×23=1702"
$tbl.Rows.Item(1).Cells.Item(2).Range.Text = "91×68=6188"
$tbl.Rows.Item(1).Cells.Item(3).Range.Text = "92×43=3956"
$tbl.Rows.Item(1).Cells.Item(4).Range.Text = "68×43=2924"
$tbl.Rows.Item(1).Cells.Item(5).Range.Text = "11×33=363"

# Row 5 holds the second line (one cell removed, shifting the rest left, plus
# a brand-new cell appended at the end).
$tbl.Rows.Item(5).Cells.Item(1).Range.Text = "53×65=3445"
$tbl.Rows.Item(5).Cells.Item(2).Range.Text = "21×90=1890"
$tbl.Rows.Item(5).Cells.Item(3).Range.Text = "95×77=7315"
$tbl.Rows.Item(5).Cells.Item(4).Range.Text = "30×18=540"
$tbl.Rows.Item(5).Cells.Item(5).Range.Text = "48×27=1296"

# Row 10 holds the third line.
$tbl.Rows.Item(10).Cells.Item(1).Range.Text = "48×67=3216"
$tbl.Rows.Item(10).Cells.Item(2).Range.Text = "29×98=2842"
$tbl.Rows.Item(10).Cells.Item(3).Range.Text = "24×98=2352"
$tbl.Rows.Item(10).Cells.Item(4).Range.Text = "90×36=3240"
$tbl.Rows.Item(10).Cells.Item(5).Range.Text = "58×86=4988"

# Row 15 holds the fourth line.
$tbl.Rows.Item(15).Cells.Item(1).Range.Text = "55×70=3850"
$tbl.Rows.Item(15).Cells.Item(2).Range.Text = "23×45=1035"
$tbl.Rows.Item(15).Cells.Item(3).Range.Text = "79×17=1343"
$tbl.Rows.Item(15).Cells.Item(4).Range.Text = "40×42=1680"
$tbl.Rows.Item(15).Cells.Item(5).Range.Text = "20×99=1980"

# Row 20 holds the fifth line.
$tbl.Rows.Item(20).Cells.Item(1).Range.Text = "93×24=2232"
$tbl.Rows.Item(20).Cells.Item(2).Range.Text = "99×76=7524"
$tbl.Rows.Item(20).Cells.Item(3).Range.Text = "73×59=4307"
$tbl.Rows.Item(20).Cells.Item(4).Range.Text = "64×37=2368"
$tbl.Rows.Item(20).Cells.Item(5).Range.Text = "33×11=363"
